# Insert a new weekly price-report row at row 3, pushing the existing
# rows 3-35 down to 4-36 (this mirrors Excel's Rows.Insert so formats,
# e.g. the date style on column D, are carried along automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Macroferia Regional de Talca"
$ws.Range("C3").Value = "Maule"
$ws.Range("D3").Value = 44532
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101001
$ws.Range("J3").Value = "Arándano (blue)"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 170
$ws.Range("N3").Value = 3600
$ws.Range("O3").Value = 3600
$ws.Range("P3").Value = 3600
$ws.Range("Q3").Value = "$/bandeja 2 kilos"
$ws.Range("R3").Value = "Provincia de Linares"
$ws.Range("S3").Value = 1800
$ws.Range("T3").Value = 2
